$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '29.017.58'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E2').Value = '  -3.91%  '
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '1.959.41'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$r = $ws.Range('D4')
$r.NumberFormat = '@'
$r.Value = '1.004'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '326.64'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  -3.65%  '
$ws.Range('E6').Value = '  +0.08%  '
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.4952'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E7').Value = '  -5.99%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.4192'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E8').Value = '  -3.88%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '52.67'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E9').Value = '  -4.12%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.09180'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E10').Value = '  -1.70%  '
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '1.094'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '22.76'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E12').Value = '  -7.03%  '
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '1.989.95'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E13').Value = '  -5.46%  '
$ws.Range('E14').Value = '  -6.00%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '7.820'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '91.24'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '0.00001093'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E18').Value = '  -5.61%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '0.06666'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E19').Value = '  -0.24%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '19.20'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  -7.92%  '
$ws.Range('E21').Value = '  +0.05%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '5.936'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  -5.77%  '
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '29.056.42'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E23').Value = '  -3.81%  '
$ws.Range('E24').Value = '  -2.93%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '2.261'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E25').Value = '  -2.37%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '2.231.75'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E26').Value = '  -3.41%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '20.54'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E27').Value = '  -5.47%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '155.42'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E28').Value = '  -4.23%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '6.256'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  -7.92%  '
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '2.245'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E30').Value = '  -9.50%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '125.99'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E31').Value = '  -5.44%  '
$ws.Range('E32').Value = '  -7.72%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '0.09807'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E33').Value = '  -6.23%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '1.520'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E34').Value = '  -8.27%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '5.827'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  -6.49%  '
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '3.679'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E36').Value = '  -5.82%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.02418'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E37').Value = '  -7.10%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '1.322'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('E39').Value = '  -8.85%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '0.06334'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  -5.60%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.6423'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E41').Value = '  -7.49%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '11.36'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  -9.15%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '0.1971'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E43').Value = '  -10.36%  '
$ws.Range('E44').Value = '  +0.14%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '0.6204'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E45').Value = '  -7.58%  '
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '1.340'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '13.36'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E47').Value = '  -6.06%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '2.187'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E48').Value = '  -7.30%  '
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '3.462'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E49').Value = '  -4.47%  '
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '0.00000000330'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -5.98%  '
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.06998'
$r.NumberFormat = 'General'
$r.Style = 'Normal'
$ws.Range('E51').Value = '  -2.96%  '
